# Fruta / hortaliza, semanal
# Insert two new weekly data rows at the top of the Pera (Packham's Triumph)
# block for "Feria Lagunitas de Puerto Montt", pushing the existing rows
# 613-637 down to 615-639.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 613 (existing rows 613:637 shift to 615:639)
$ws.Rows("613:614").Insert()

# New row 613: Pera, Packham's Triumph, Primera
$ws.Cells.Item(613, 1).Value = 4
$ws.Cells.Item(613, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(613, 3).Value = "Los Lagos"
$ws.Cells.Item(613, 4).Value = 45267
$ws.Cells.Item(613, 5).Value = 10
$ws.Cells.Item(613, 6).Value = "Fruta"
$ws.Cells.Item(613, 7).Value = 100104
$ws.Cells.Item(613, 8).Value = "Frutos de pepita"
$ws.Cells.Item(613, 9).Value = 100104005
$ws.Cells.Item(613, 10).Value = "Pera"
$ws.Cells.Item(613, 11).Value = "Packham's Triumph"
$ws.Cells.Item(613, 12).Value = "Primera"
$ws.Cells.Item(613, 13).Value = 250
$ws.Cells.Item(613, 14).Value = 21000
$ws.Cells.Item(613, 15).Value = 21000
$ws.Cells.Item(613, 16).Value = 21000
$ws.Cells.Item(613, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(613, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(613, 19).Value = 1400
$ws.Cells.Item(613, 20).Value = 15

# New row 614: Pera, Packham's Triumph, Segunda
$ws.Cells.Item(614, 1).Value = 4
$ws.Cells.Item(614, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(614, 3).Value = "Los Lagos"
$ws.Cells.Item(614, 4).Value = 45267
$ws.Cells.Item(614, 5).Value = 10
$ws.Cells.Item(614, 6).Value = "Fruta"
$ws.Cells.Item(614, 7).Value = 100104
$ws.Cells.Item(614, 8).Value = "Frutos de pepita"
$ws.Cells.Item(614, 9).Value = 100104005
$ws.Cells.Item(614, 10).Value = "Pera"
$ws.Cells.Item(614, 11).Value = "Packham's Triumph"
$ws.Cells.Item(614, 12).Value = "Segunda"
$ws.Cells.Item(614, 13).Value = 250
$ws.Cells.Item(614, 14).Value = 16000
$ws.Cells.Item(614, 15).Value = 16000
$ws.Cells.Item(614, 16).Value = 16000
$ws.Cells.Item(614, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(614, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(614, 19).Value = 1067
$ws.Cells.Item(614, 20).Value = 15
